# Word COM-interop script: replace the placeholder "Fiche de contact" card
# with six anonymized-style contact paragraphs, each containing three
# line-break-separated sentences. Net paragraph count: 8 -> 6.
$d = $word.ActiveDocument

function Set-ParaContent($para, $sentences) {
    $inner = ""
    foreach ($s in $sentences) {
        $inner += "<w:t>" + $s + "</w:t><w:br/>"
    }
    $xml = "<w:p xmlns:w=`"http://schemas.openxmlformats.org/wordprocessingml/2006/main`"><w:r>" + $inner + "</w:r></w:p>"
    $para.Range.InsertXML($xml)
}

$sentences1 = @(
    "Célina Gosselin (annemarechal@example.org), né le 27/07/1967 à SanchezVille, travaille chez Cordier Gaillard S.A.R.L..",
    "Claire Guillou (gillesgregoire@example.com), né le 19/12/1952 à Da Silva-sur-Mendès, travaille chez Bernier SARL.",
    "Ils participeront à un événement le 1979-04-17 à Ruiz, organisé par Giraud Payet S.A.S..",
)
Set-ParaContent $d.Paragraphs(1) $sentences1

$sentences2 = @(
    "Thomas Michel-Thibault (gguillet@example.com), né le 16/05/1953 à Saint Margot-la-Forêt, travaille chez Roy Gros S.A.R.L..",
    "Luce Lévy (genevievecourtois@example.com), né le 07/06/1916 à Saint Denis-sur-Mer, travaille chez Vallet.",
    "Ils participeront à un événement le 1981-01-08 à Reynaud, organisé par Bailly.",
)
Set-ParaContent $d.Paragraphs(2) $sentences2

$sentences3 = @(
    "Xavier Lecomte (leonmaillet@example.net), né le 31/12/2023 à Besnard-sur-Mer, travaille chez Valentin.",
    "Nathalie Evrard (delattrejosephine@example.com), né le 04/09/1954 à Blanchard-sur-Langlois, travaille chez Samson Mercier S.A.S..",
    "Ils participeront à un événement le 1974-12-30 à De Sousa-sur-Mer, organisé par Clément Bonnet SA.",
)
Set-ParaContent $d.Paragraphs(3) $sentences3

$sentences4 = @(
    "Bernard-Daniel Seguin (camille79@example.org), né le 10/10/1973 à Dupont-les-Bains, travaille chez Devaux.",
    "Georges-Robert Potier (simonemarchand@example.com), né le 31/10/2001 à Da Silva-la-Forêt, travaille chez Martineau Ramos S.A.R.L..",
    "Ils participeront à un événement le 1970-09-04 à Costa, organisé par Hamel Renard S.A.S..",
)
Set-ParaContent $d.Paragraphs(4) $sentences4

$sentences5 = @(
    "Clémence Lévêque (aweber@example.com), né le 15/08/2004 à Meunier, travaille chez Chauveau Joly SA.",
    "Gabriel Vincent (nicolegilbert@example.net), né le 20/01/2023 à Barre, travaille chez Albert.",
    "Ils participeront à un événement le 1976-01-17 à Saint Thierry, organisé par Benoit.",
)
Set-ParaContent $d.Paragraphs(5) $sentences5

$sentences6 = @(
    "Aimé Gimenez (carreagnes@example.net), né le 29/08/1920 à Sainte Maurice, travaille chez Grégoire SARL.",
    "Gabriel Blanc-Guillou (sabine94@example.org), né le 17/12/1927 à LévêqueVille, travaille chez Joly.",
    "Ils participeront à un événement le 1987-08-06 à Julien-sur-Fournier, organisé par Picard SARL.",
)
Set-ParaContent $d.Paragraphs(6) $sentences6

# The original document had 8 paragraphs (1 heading + 7 fields); the new
# content only needs 6, so drop the last two (formerly "Lieu de naissance"
# and "Employeur") now that their text has been folded into paragraph 6.
$d.Paragraphs(8).Range.Delete()
$d.Paragraphs(7).Range.Delete()

Write-Output ("Final paragraph count: " + $d.Paragraphs.Count)